# feat: add 2022-Q1 data
#
# Plan:
# 1. Duplicate the existing "总计" sheet and place the copy right after "2021-Q4".
#    This copy keeps the original sheet's sheetPr/pageMargins/styles intact, and
#    will become the new "2022-Q1" sheet.
# 2. Rename the copy to "2022-Q1", extend it to 8 columns / 7 rows (re-using the
#    existing header/A-column style for the new cells) and fill in the fund
#    holdings data for 2022-Q1.
# 3. On the original "总计" sheet, insert a new row 2 (copying formats from a
#    stable, unrelated cell so the inserted row ends up using the correct
#    existing style rather than creating new ones) and fill it in with the
#    2022-Q1 summary figures, renumbering the index column below it.

$wb = $excel.ActiveWorkbook

$zj = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# Step 1 + 2: build the new "2022-Q1" sheet from a copy of "总计"
# ---------------------------------------------------------------------------
$zj.Copy($null, $q4)
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Extend the header style (from D1) across the new columns E1:H1
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Extend the "A" index-column style (from row 3) down across the new rows 4:7
$q1.Range("A3:D3").Copy()
$q1.Range("A4:D7").PasteSpecial(-4122)

# Extend the plain data-cell style (from B2:D2) across the new columns/rows
$q1.Range("B2:D2").Copy()
$q1.Range("E2:H7").PasteSpecial(-4122)

# Headers
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund holdings data for 2022-Q1 (columns B-G are stored as text, like the
# other quarterly sheets; column H holds a numeric rank)
$fundRows = @(
  @("001822", "华商智能生活灵活配置混合", "21.89", "90.29", "8.11", "1.7753", 2),
  @("009693", "富国积极成长一年定期开放混合", "17.82", "98.74", "3.37", "0.6005", 6),
  @("001933", "华商新兴活力灵活配置混合", "6.91", "92.28", "7.95", "0.5493", 1),
  @("001521", "国寿安保成长优选股票", "4.24", "87.83", "4.18", "0.1772", 7),
  @("008082", "国寿安保研究精选混合A", "0.52", "91.60", "4.80", "0.0250", 7),
  @("008083", "国寿安保研究精选混合C", "0.15", "91.60", "4.80", "0.0072", 7)
)

# Force text storage for the B:G data block (preserves leading zeros / trailing
# zeros in codes and figures) and then strip the format back off again so the
# cells end up with no explicit style, matching the rest of the sheet.
$dataBlock = $q1.Range("B2:G7")
$dataBlock.NumberFormat = "@"

$r = 2
foreach ($row in $fundRows) {
    $q1.Range("A$r").Value = $r - 2
    $q1.Range("B$r").Value = $row[0]
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").Value = $row[2]
    $q1.Range("E$r").Value = $row[3]
    $q1.Range("F$r").Value = $row[4]
    $q1.Range("G$r").Value = $row[5]
    $q1.Range("H$r").Value = $row[6]
    $r++
}

$dataBlock.ClearFormats()

# ---------------------------------------------------------------------------
# Step 3: insert the new summary row into the original "总计" sheet
# ---------------------------------------------------------------------------
# NOTE: the original $zj reference became stale once it was Copy()'d above
# (it now points at the copy instead of the original sheet), so look the
# "总计" sheet up again fresh by name before editing it.
$zj = $wb.Worksheets.Item("总计")

$zj.Rows.Item(2).Insert()

# The row-insert carries the formatting of the row above down into the new
# row; reset it and instead copy the correct existing style (from the
# untouched "2021-Q4" sheet's A3 cell, which isn't affected by the insert)
# onto the new index cell.
$zj.Range("B2:D2").ClearFormats()
$q4.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 6
$zj.Range("D2").Value = 3.13

$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
